$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cluster Name"
$ws.Range("B1").Value = "Active Cases"

$ws.Range("A2").Value = "12 Ironbark Drive Sunbury"
$ws.Range("B2").Value = 12
$ws.Range("A3").Value = "202111 45784 Holy Rosary Primary School White Hills"
$ws.Range("B3").Value = 36
$ws.Range("A4").Value = "3321 Rochester and Elmore District Health Service Yalunkan Aged Care Hostel Rochester"
$ws.Range("B4").Value = 15
$ws.Range("A5").Value = "3323 Villa Maria Catholic Homes St Bernadette's Aged Care Sunshine North"
$ws.Range("B5").Value = 13
$ws.Range("A6").Value = "3600 Belvedere Aged Care Noble Park"
$ws.Range("B6").Value = 22
$ws.Range("A7").Value = "3601 Baptcare Westhaven community"
$ws.Range("B7").Value = 25
$ws.Range("A8").Value = "3653 Fronditha Thalpori St Albans Aged Care"
$ws.Range("B8").Value = 26
$ws.Range("A9").Value = "4257 BlueCross The Gables Camberwell"
$ws.Range("B9").Value = 27
$ws.Range("A10").Value = "44098 Stawell Primary School"
$ws.Range("B10").Value = 22
$ws.Range("A11").Value = "44121 Wallan Primary School Wallan"
$ws.Range("B11").Value = 20
$ws.Range("A12").Value = "44165 Greenvale Primary School"
$ws.Range("B12").Value = 25
$ws.Range("A13").Value = "44234 Lucknow Primary School Bairnsdale"
$ws.Range("B13").Value = 26
$ws.Range("A14").Value = "44495 Lakes Entrance Primary School"
$ws.Range("B14").Value = 10
$ws.Range("A15").Value = "44667 Beaumaris Primary School Beaumaris"
$ws.Range("B15").Value = 22
$ws.Range("A16").Value = "44701 Hampton Park Primary School Hampton Park"
$ws.Range("B16").Value = 11
$ws.Range("A17").Value = "44718 Parkdale Primary School Parkdale"
$ws.Range("B17").Value = 12
$ws.Range("A18").Value = "44811 Dandenong North Primary School Dandenong"
$ws.Range("B18").Value = 34
$ws.Range("A19").Value = "44865 Parktone Primary School Parkdale"
$ws.Range("B19").Value = 28
$ws.Range("A20").Value = "44891 Cranbourne Park Primary School Cranbourne"
$ws.Range("B20").Value = 18
$ws.Range("A21").Value = "45248 Brookside P-9 College Caroline Springs"
$ws.Range("B21").Value = 15
$ws.Range("A22").Value = "45249 Creekside K-9 College Caroline Springs"
$ws.Range("B22").Value = 16
$ws.Range("A23").Value = "45267 Epping Views Primary School Epping"
$ws.Range("B23").Value = 12
$ws.Range("A24").Value = "45569 Nhill College Nhill"
$ws.Range("B24").Value = 34
$ws.Range("A25").Value = "45648 St Brendans Primary School Shepparton"
$ws.Range("B25").Value = 14
$ws.Range("A26").Value = "4574 Village Glen Aged Care Residences Mornington"
$ws.Range("B26").Value = 15
$ws.Range("A27").Value = "45836 St Joseph's Primary School Sorrento"
$ws.Range("B27").Value = 12
$ws.Range("A28").Value = "45967 St Clement of Rome School Bulleen"
$ws.Range("B28").Value = 10
$ws.Range("A29").Value = "46037 Nazareth Catholic Primary School Grovedale"
$ws.Range("B29").Value = 35
$ws.Range("A30").Value = "46050 Our Lady's Catholic Primary School Craigieburn"
$ws.Range("B30").Value = 21
$ws.Range("A31").Value = "46125 Our Lady of the Southern Cross Primary School Manor Lakes"
$ws.Range("B31").Value = 20
$ws.Range("A32").Value = "46190 Haileybury Brighton East"
$ws.Range("B32").Value = 11
$ws.Range("A33").Value = "46215 Yeshivah Primary College St Kilda East"
$ws.Range("B33").Value = 10
$ws.Range("A34").Value = "46276 Hillcrest Christian College Clyde North"
$ws.Range("B34").Value = 18
$ws.Range("A35").Value = "46328 Ilim College Dallas Main Campus Dallas"
$ws.Range("B35").Value = 24
$ws.Range("A36").Value = "50395 St Francis of Assisi Catholic PrimarySchool Tarneit"
$ws.Range("B36").Value = 10
$ws.Range("A37").Value = "52380 Al Iman College Melton South"
$ws.Range("B37").Value = 11
$ws.Range("A38").Value = "52473 John Henry Primary School Pakenham"
$ws.Range("B38").Value = 19
$ws.Range("A39").Value = "Camp Coolamatong Farm Camp Banksia Peninsula"
$ws.Range("B39").Value = 12
$ws.Range("A40").Value = "Escala NewQuay Construction Site Docklands Drive Docklands"
$ws.Range("B40").Value = 16
$ws.Range("A41").Value = "Gippsland and East Gippsland Aboriginal Co-Operative Bairnsdale"
$ws.Range("B41").Value = 13
$ws.Range("A42").Value = "Hamilton Country Music Festival Hamilton Golf Club Hamilton"
$ws.Range("B42").Value = 29
$ws.Range("A43").Value = "Melton Willows Melton"
$ws.Range("B43").Value = 10
$ws.Range("A44").Value = "St Josephs Catholic Primary School Warragul"
$ws.Range("B44").Value = 11
$ws.Range("A45").Value = "Wagstaff Meat Processing Plant Cranbourne East"
$ws.Range("B45").Value = 23
